$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value. The Price/Volume columns contain
# strings that look numeric (e.g. "250.48") or percent-like (e.g. "  -0.96%  ").
# Excel would otherwise auto-convert these into numbers on assignment, so we
# briefly force the cell to Text format, write the value, then restore the
# original General format so the cell keeps behaving like the rest of the sheet.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

Set-TextValue "D2" "41.375.73"
Set-TextValue "E2" "  -0.96%  "
Set-TextValue "D3" "2.189.78"
Set-TextValue "E3" "  -1.23%  "
Set-TextValue "D5" "250.48"
Set-TextValue "E5" "  -0.15%  "
Set-TextValue "D6" "0.616"
Set-TextValue "E6" "  -2.20%  "
Set-TextValue "D7" "67.43"
Set-TextValue "E7" "  -4.60%  "
Set-TextValue "D9" "0.593"
Set-TextValue "E9" "  -1.83%  "
Set-TextValue "D10" "39.17"
Set-TextValue "E10" "  -2.27%  "
Set-TextValue "D11" "59.74"
Set-TextValue "E11" "  +2.58%  "
Set-TextValue "D12" "0.0943"
Set-TextValue "E12" "  -2.23%  "
Set-TextValue "E13" "  -0.93%  "
Set-TextValue "D14" "6.95"
Set-TextValue "E14" "  -3.87%  "
Set-TextValue "D15" "2.515.11"
Set-TextValue "E15" "  -1.25%  "
Set-TextValue "D16" "14.48"
Set-TextValue "E16" "  -3.26%  "
Set-TextValue "D17" "0.855"
Set-TextValue "E17" "  -3.89%  "
Set-TextValue "D18" "2.198.21"
Set-TextValue "E18" "  -0.95%  "
Set-TextValue "D19" "41.319.80"
Set-TextValue "E19" "  -1.15%  "
Set-TextValue "D20" "0.0₃0950"
Set-TextValue "E20" "  -1.14%  "
Set-TextValue "D21" "71.92"
Set-TextValue "E21" "  -0.90%  "
Set-TextValue "D22" "6.10"
Set-TextValue "E22" "  -2.34%  "
Set-TextValue "D23" "230.86"
Set-TextValue "E23" "  -1.66%  "
Set-TextValue "D24" "2.04"
Set-TextValue "E24" "  -1.73%  "
Set-TextValue "D25" "3.84"
Set-TextValue "E25" "  -6.26%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D26" "11.44"
Set-TextValue "E26" "  -2.66%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  +0.14%  "
Set-TextValue "E28" "  -4.24%  "
Set-TextValue "E29" "  -2.16%  "
Set-TextValue "E30" "  -2.86%  "
Set-TextValue "D31" "167.24"
Set-TextValue "E31" "  -1.61%  "
Set-TextValue "D32" "20.25"
Set-TextValue "E32" "  -2.59%  "
Set-TextValue "E33" "  -1.36%  "
Set-TextValue "D34" "5.72"
Set-TextValue "E34" "  +3.57%  "
Set-TextValue "D35" "0.0764"
Set-TextValue "E35" "  +3.70%  "
Set-TextValue "E36" "  -1.93%  "
Set-TextValue "D37" "4.21"
Set-TextValue "E37" "  +3.97%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D38" "25.91"
Set-TextValue "E38" "  -0.52%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D39" "4.55"
Set-TextValue "E39" "  -2.53%  "
Set-TextValue "D40" "0.0304"
Set-TextValue "E40" "  -0.91%  "
Set-TextValue "E41" "  -2.30%  "
Set-TextValue "D42" "5.20"
Set-TextValue "E42" "  +8.47%  "
Set-TextValue "E43" "  -4.93%  "
Set-TextValue "D44" "11.63"
Set-TextValue "E44" "  -6.46%  "
Set-TextValue "D45" "61.16"
Set-TextValue "E45" "  -5.80%  "
Set-TextValue "E46" "  -7.05%  "
Set-TextValue "E47" "  -1.87%  "
Set-TextValue "D48" "0.0991"
Set-TextValue "E48" "  -2.58%  "
Set-TextValue "E49" "  -0.04%  "
Set-TextValue "D50" "1.15"
Set-TextValue "E50" "  -1.46%  "
Set-TextValue "D51" "4.35"
Set-TextValue "E51" "  -6.58%  "
